$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3528.4285
$ws.Range("I40").Value = 2540
$ws.Range("K40").Value = 2540
$ws.Range("M40").Value = -2365
$ws.Range("H76").Value = 4400.048
$ws.Range("I76").Value = 3600.5715
$ws.Range("J76").Value = 5999
$ws.Range("K76").Value = 3600.5715
$ws.Range("L76").Value = 5999
$ws.Range("M76").Value = -3285.5715
$ws.Range("N76").Value = -6629
$ws.Range("H79").Value = 4400.048
$ws.Range("I79").Value = 3600.5715
$ws.Range("J79").Value = 5999
$ws.Range("K79").Value = 3600.5715
$ws.Range("L79").Value = 5999
$ws.Range("M79").Value = -2508.5715
$ws.Range("N79").Value = -8183
$ws.Range("H112").Value = 6142.815
$ws.Range("I112").Value = 662.25
$ws.Range("K112").Value = 1986.75
$ws.Range("M112").Value = -878.75
$ws.Range("H113").Value = 4771.143
$ws.Range("J113").Value = 5798.6
$ws.Range("L113").Value = 5798.6
$ws.Range("N113").Value = -12306.6
$ws.Range("H137").Value = 12928.821
$ws.Range("I137").Value = 16243.904
$ws.Range("K137").Value = 48731.712
$ws.Range("M137").Value = -46181.712
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 349
$ws.Range("I4").Value = 298
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 298
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -182
$ws.Range("N4").Value = -632
$ws.Range("H5").Value = 7540.4287
$ws.Range("I5").Value = 9529.182000000001
$ws.Range("K5").Value = 9529.182000000001
$ws.Range("M5").Value = -9417.182000000001
$ws.Range("H50").Value = 1369
$ws.Range("I50").Value = 48
$ws.Range("K50").Value = 48
$ws.Range("M50").Value = 666
$ws.Range("H61").Value = 4168.9688
$ws.Range("I61").Value = 800.2857
$ws.Range("K61").Value = 800.2857
$ws.Range("M61").Value = -588.2857
$ws.Range("H122").Value = 2170.8572
$ws.Range("I122").Value = 2078.2727
$ws.Range("K122").Value = 6234.8181
$ws.Range("M122").Value = -3784.8181
$ws.Range("H136").Value = 4168.9688
$ws.Range("I136").Value = 800.2857
$ws.Range("K136").Value = 2400.8571
$ws.Range("M136").Value = 149.1428999999998
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 7540.4287
$ws.Range("I4").Value = 9529.182000000001
$ws.Range("K4").Value = 9529.182000000001
$ws.Range("M4").Value = -9414.182000000001
$ws.Range("H20").Value = 13123.08
$ws.Range("I20").Value = 17830.277
$ws.Range("J20").Value = 1018.8571
$ws.Range("K20").Value = 17830.277
$ws.Range("L20").Value = 1018.8571
$ws.Range("M20").Value = -17583.277
$ws.Range("N20").Value = -1512.8571
$ws.Range("H22").Value = 1001
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H86").Value = 9633
$ws.Range("I86").Value = 9000
$ws.Range("K86").Value = 9000
$ws.Range("M86").Value = -7877
$ws.Range("H89").Value = 9633
$ws.Range("I89").Value = 9000
$ws.Range("K89").Value = 45000
$ws.Range("M89").Value = -39384
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 14333
$ws.Range("I17").Value = 14333
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 14333
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -14159
$ws.Range("H31").Value = 5265696.5
$ws.Range("I31").Value = 11112110
$ws.Range("K31").Value = 11112110
$ws.Range("M31").Value = -11111815
$ws.Range("H34").Value = 5265696.5
$ws.Range("I34").Value = 11112110
$ws.Range("K34").Value = 11112110
$ws.Range("M34").Value = -11111908
$ws.Range("H86").Value = 50381
$ws.Range("I86").Value = 60450.77
$ws.Range("K86").Value = 60450.77
$ws.Range("M86").Value = -59327.77
$ws.Range("H89").Value = 50381
$ws.Range("I89").Value = 60450.77
$ws.Range("K89").Value = 302253.85
$ws.Range("M89").Value = -296637.85
$ws.Range("H99").Value = 13165.667
$ws.Range("I99").Value = 11661.667
$ws.Range("J99").Value = 14669.667
$ws.Range("K99").Value = 11661.667
$ws.Range("L99").Value = 14669.667
$ws.Range("M99").Value = -10163.667
$ws.Range("N99").Value = -17665.667
$ws.Range("H126").Value = 13165.667
$ws.Range("I126").Value = 11661.667
$ws.Range("J126").Value = 14669.667
$ws.Range("K126").Value = 34985.001
$ws.Range("L126").Value = 44009.001
$ws.Range("M126").Value = -32515.001
$ws.Range("N126").Value = -48949.001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 240.26666
$ws.Range("I2").Value = 270.5
$ws.Range("K2").Value = 1623
$ws.Range("M2").Value = -1510
$ws.Range("H120").Value = 15871.5
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = ""
$ws.Range("H132").Value = 1180.6364
$ws.Range("I132").Value = 851
$ws.Range("K132").Value = 7659
$ws.Range("M132").Value = -5129
$ws.Range("H138").Value = 5544.0835
$ws.Range("I138").Value = 1166.2222
$ws.Range("J138").Value = 18677.666
$ws.Range("K138").Value = 3498.6666
$ws.Range("L138").Value = 56032.99800000001
$ws.Range("M138").Value = 1641.3334
$ws.Range("N138").Value = -66312.99800000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 689.7778
$ws.Range("I2").Value = 874
$ws.Range("J2").Value = 321.33334
$ws.Range("K2").Value = 874
$ws.Range("L2").Value = 321.33334
$ws.Range("M2").Value = -761
$ws.Range("N2").Value = -547.33334
$ws.Range("H35").Value = 1000000
$ws.Range("I35").Value = 1000000
$ws.Range("K35").Value = 1000000
$ws.Range("M35").Value = -999702
$ws.Range("H80").Value = 8776.333000000001
$ws.Range("J80").Value = 17249.5
$ws.Range("L80").Value = 17249.5
$ws.Range("N80").Value = -19245.5
$ws.Range("H83").Value = 8776.333000000001
$ws.Range("J83").Value = 17249.5
$ws.Range("L83").Value = 86247.5
$ws.Range("N83").Value = -96231.5
$ws.Range("H92").Value = 200
$ws.Range("J92").Value = 200
$ws.Range("L92").Value = 200
$ws.Range("N92").Value = -3944
$ws.Range("H132").Value = 2159.1765
$ws.Range("I132").Value = 1477.7693
$ws.Range("K132").Value = 4433.3079
$ws.Range("M132").Value = -1903.3079
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 39500
$ws.Range("J11").Value = 39500
$ws.Range("L11").Value = 39500
$ws.Range("N11").Value = -39780
$ws.Range("H13").Value = 6005
$ws.Range("J13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("N13").Value = -12280
$ws.Range("H82").Value = 1276.5625
$ws.Range("I82").Value = 1103.1111
$ws.Range("K82").Value = 1103.1111
$ws.Range("M82").Value = -742.1111000000001
$ws.Range("H85").Value = 1276.5625
$ws.Range("I85").Value = 1103.1111
$ws.Range("K85").Value = 1103.1111
$ws.Range("M85").Value = 144.8888999999999
$ws.Range("H100").Value = 3314
$ws.Range("I100").Value = 2974.5
$ws.Range("J100").Value = 3766.6667
$ws.Range("K100").Value = 2974.5
$ws.Range("L100").Value = 3766.6667
$ws.Range("M100").Value = -2433.5
$ws.Range("N100").Value = -4848.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 10000000
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = ""
$ws.Range("H122").Value = 38236.6
$ws.Range("I122").Value = 45469.332
$ws.Range("J122").Value = 4139.4287
$ws.Range("K122").Value = 136407.996
$ws.Range("L122").Value = 12418.2861
$ws.Range("M122").Value = -133957.996
$ws.Range("N122").Value = -17318.2861
$ws.Range("H136").Value = 17203.027
$ws.Range("I136").Value = 18576.182
$ws.Range("K136").Value = 55728.546
$ws.Range("M136").Value = -53178.546
